$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.06030000000001
$ws.Range("B3").Value = 6.038500000000002
$ws.Range("D3").Value = -7.125899999999998
$ws.Range("E6").Value = 16.4328
$ws.Range("D12").Value = -7.316899999999999
$ws.Range("B14").Value = 5.619800000000001
$ws.Range("E19").Value = 16.3291
$ws.Range("B21").Value = 9.511900000000006
$ws.Range("B23").Value = 9.139400000000004
$ws.Range("D24").Value = -7.360000000000005
$ws.Range("E24").Value = 16.66490000000001
$ws.Range("B25").Value = 5.358900000000003
$ws.Range("C25").Value = -13.1084
$ws.Range("D25").Value = -9.084299999999992
$ws.Range("B26").Value = 5.547300000000003
$ws.Range("C27").Value = -12.9893
$ws.Range("B29").Value = 5.027500000000002
$ws.Range("E30").Value = 15.33929999999999
$ws.Range("C31").Value = -13.19619999999999
$ws.Range("E31").Value = 16.02320000000001
$ws.Range("E33").Value = 17.06100000000002
$ws.Range("C39").Value = -12.69260000000001
$ws.Range("E42").Value = 16.55940000000001
$ws.Range("C48").Value = -11.50379999999999
$ws.Range("D50").Value = -7.990000000000001
$ws.Range("C51").Value = -11.56229999999999
$ws.Range("C52").Value = -11.1372
$ws.Range("B53").Value = 5.347299999999998
$ws.Range("D53").Value = -6.047700000000001
$ws.Range("C55").Value = -13.7109
$ws.Range("E55").Value = 16.3392
$ws.Range("C56").Value = -12.07769999999999
$ws.Range("B57").Value = 5.014299999999999
$ws.Range("C57").Value = -13.736
$ws.Range("D57").Value = -8.841099999999999
$ws.Range("E58").Value = 16.47640000000002
$ws.Range("B59").Value = 4.870799999999999
$ws.Range("D61").Value = -7.803399999999999
$ws.Range("D63").Value = -7.692700000000002
$ws.Range("E65").Value = 17.10690000000001
$ws.Range("B69").Value = 5.415899999999993
$ws.Range("D70").Value = -8.1495
$ws.Range("E70").Value = 16.65969999999999
$ws.Range("C73").Value = -12.437
$ws.Range("E75").Value = 16.57920000000001
$ws.Range("B79").Value = 9.465200000000005
$ws.Range("B83").Value = 5.406799999999999
$ws.Range("E83").Value = 16.6621
$ws.Range("D86").Value = -8.688899999999997
$ws.Range("E86").Value = 16.07170000000001
$ws.Range("C89").Value = -10.8032
$ws.Range("C90").Value = -12.9025
$ws.Range("B91").Value = 5.062299999999998
$ws.Range("C92").Value = -11.3428
$ws.Range("B93").Value = 5.941999999999998
$ws.Range("E96").Value = 15.8358
$ws.Range("E97").Value = 16.85340000000001
$ws.Range("D98").Value = -9.185599999999992
$ws.Range("D100").Value = -8.861699999999997
$ws.Range("D102").Value = -7.982799999999997
